$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data adds two new weekly price records (rows 55-56) for
# "Femacal de La Calera - Tuna" dated 2022-03-09 (serial 44629), pushing the
# existing rows 55-115 down to 57-117.
$ws.Rows("55:56").Insert()

# New row 55: Primera, 45 units, $15000, $/caja 16 kilos, Cabildo, $938/kg, 16 kg/unit
$ws.Range("A55").Value = 3
$ws.Range("B55").Value = "Femacal de La Calera"
$ws.Range("C55").Value = "Coquimbo"
$ws.Range("D55").Value = 44629
$ws.Range("E55").Value = 5
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100107
$ws.Range("H55").Value = "Otros"
$ws.Range("I55").Value = 100107011
$ws.Range("J55").Value = "Tuna"
$ws.Range("K55").Value = "Sin especificar"
$ws.Range("L55").Value = "Primera"
$ws.Range("M55").Value = 45
$ws.Range("N55").Value = 15000
$ws.Range("O55").Value = 15000
$ws.Range("P55").Value = 15000
$ws.Range("Q55").Value = "$/caja 16 kilos"
$ws.Range("R55").Value = "Cabildo"
$ws.Range("S55").Value = 938
$ws.Range("T55").Value = 16

# New row 56: Segunda, 50 units, $13000, $/caja 16 kilos, Cabildo, $812/kg, 16 kg/unit
$ws.Range("A56").Value = 3
$ws.Range("B56").Value = "Femacal de La Calera"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 44629
$ws.Range("E56").Value = 5
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100107
$ws.Range("H56").Value = "Otros"
$ws.Range("I56").Value = 100107011
$ws.Range("J56").Value = "Tuna"
$ws.Range("K56").Value = "Sin especificar"
$ws.Range("L56").Value = "Segunda"
$ws.Range("M56").Value = 50
$ws.Range("N56").Value = 13000
$ws.Range("O56").Value = 13000
$ws.Range("P56").Value = 13000
$ws.Range("Q56").Value = "$/caja 16 kilos"
$ws.Range("R56").Value = "Cabildo"
$ws.Range("S56").Value = 812
$ws.Range("T56").Value = 16
